$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D (Price) and E (Volume(1h)) keep their values stored as text,
# matching the workbook author's convention of inline-string numeric/percentage text.
$ws.Columns("D:E").NumberFormat = "@"

# Apply the updated symbol list values cell by cell
$ws.Range('D2').Value = '303.08'
$ws.Range('E2').Value = '2.67%'
$ws.Range('D3').Value = '35.05'
$ws.Range('E3').Value = '12.66%'
$ws.Range('D4').Value = '5.153'
$ws.Range('E4').Value = '4.35%'
$ws.Range('D5').Value = '0.07776'
$ws.Range('E5').Value = '4.58%'
$ws.Range('D6').Value = '2.372'
$ws.Range('E6').Value = '6.31%'
$ws.Range('D7').Value = '8.033'
$ws.Range('E7').Value = '3.45%'
$ws.Range('D8').Value = '3.937'
$ws.Range('E8').Value = '5.25%'
$ws.Range('D9').Value = '0.9326'
$ws.Range('E9').Value = '1.98%'
$ws.Range('D10').Value = '0.09947'
$ws.Range('E10').Value = '10.03%'
$ws.Range('D11').Value = '0.1797'
$ws.Range('E11').Value = '4.85%'
$ws.Range('D12').Value = '0.08633'
$ws.Range('E12').Value = '4.51%'
$ws.Range('D13').Value = '0.03318'
$ws.Range('E13').Value = '5.58%'
$ws.Range('D14').Value = '0.09930'
$ws.Range('E14').Value = '-1.45%'
$ws.Range('D15').Value = '0.001497'
$ws.Range('E15').Value = '-1.05%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '0.005770'
$ws.Range('E16').Value = '0.89%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.461'
$ws.Range('E17').Value = '-1.06%'
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').Value = '2.142'
$ws.Range('E18').Value = '3.01%'
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').Value = '0.3359'
$ws.Range('E19').Value = '0.94%'
$ws.Range('B20').Value = 'ProBitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D20').Value = '0.1333'
$ws.Range('E20').Value = '3.09%'
$ws.Range('B21').Value = 'MCDex'
$ws.Range('C21').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D21').Value = '4.286'
$ws.Range('E21').Value = '7.38%'
$ws.Range('B22').Value = 'ZBToken'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D22').Value = '0.2299'
$ws.Range('E22').Value = '9.44%'
$ws.Range('B23').Value = 'CoinExToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D23').Value = '0.04562'
$ws.Range('E23').Value = '0.26%'
$ws.Range('E24').Value = '-0.07%'
$ws.Range('D25').Value = '0.004369'
$ws.Range('E25').Value = '-5.32%'
$ws.Range('E26').Value = '-0.15%'
$ws.Range('E27').Value = '-0.02%'
$ws.Range('D39').Value = '0.01789'
$ws.Range('E39').Value = '10.91%'
$ws.Range('D40').Value = '0.04796'
$ws.Range('E40').Value = '5.95%'
$ws.Range('D41').Value = '0.007749'
$ws.Range('E41').Value = '6.86%'
$ws.Range('E42').Value = '6.16%'
$ws.Range('D43').Value = '0.007103'
$ws.Range('E43').Value = '-21.12%'
$ws.Range('E44').Value = '9.23%'
$ws.Range('D45').Value = '0.009440'
$ws.Range('E45').Value = '3.30%'
$ws.Range('E46').Value = '0.38%'
$ws.Range('E47').Value = '-0.18%'
$ws.Range('D48').Value = '3.062'
$ws.Range('E48').Value = '33.84%'
$ws.Range('E49').Value = '-0.09%'
$ws.Range('E50').Value = '-0.18%'
$ws.Range('E51').Value = '-0.18%'
